# Restore revision #7d3c5372341a2951d88b1ca9233dc3dc42ac2095.TEST
#
# Sheet "Rules":
#   - C10 (numeric literal): 18 -> 1
#   - B11 (string literal):  "1" -> "R40"
# Once "1" is no longer referenced by any cell, Excel drops it from the
# shared-strings table automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1.0
$ws.Range("B11").Value = "R40"
